$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows for 2021/11/26 (date serial 44526) appended below the
# existing table (previously ending at row 20).
# Columns: A=Date, B=Stock_Id, C=ProfitPercent, D=ProfitMoney
$newRows = @(
    @(44526, 2340, 0.1, 290.6999999999998),
    @(44526, 6104, 0.05, 144),
    @(44526, 6138, 0.25, 765),
    @(44526, 3016, 0.01, 21),
    @(44526, 8069, 0.03, 87.5),
    @(44526, 8289, 0.04, 101.2499999999995),
    @(44526, 1712, 0.03, 80.25),
    @(44526, 3588, -0.02, -87.5),
    @(44526, 5351, -0.07000000000000001, -403.2000000000007),
    @(44526, 2484, -0.06, -349.3500000000004)
)

$startRow = 21
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $row[0]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Columns E-K stay blank (as in the rest of the table) but are written
    # as empty-text cells rather than left completely absent, matching
    # the empty Action/PositionSize/Price/EachCost columns used elsewhere.
    for ($c = 5; $c -le 11; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
}
